$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Mark values for questions 1-3 ---
$ws.Range("I2").Value = 4
$ws.Range("I3").Value = 2
$ws.Range("I4").Value = 4

# --- Remove question 4 ("Django is written in C++?") content, leaving row 5 blank ---
$ws.Range("A5:I5").ClearContents()

# Row 5 no longer carries the table border - blank spacer row
$ws.Range("A5:I5").Borders.LineStyle = 0

# A5:F5 and I5 keep the "answer" look (Arial 9, vertical-centered, wrapped) minus the border
$ws.Range("A5:F5").Font.Name = "Arial"
$ws.Range("A5:F5").Font.Size = 9
$ws.Range("A5:F5").VerticalAlignment = -4108
$ws.Range("A5:F5").WrapText = $true

$ws.Range("I5").Font.Name = "Arial"
$ws.Range("I5").Font.Size = 9
$ws.Range("I5").VerticalAlignment = -4108
$ws.Range("I5").WrapText = $true

# G5:H5 keep the "blank option cell" look (default font, top-aligned, wrapped) minus the border
$ws.Range("G5:H5").VerticalAlignment = -4160
$ws.Range("G5:H5").WrapText = $true

# --- Update the active selection to reflect where the user ended up ---
$ws.Range("J6").Select() | Out-Null
